$wb = $excel.ActiveWorkbook

# --- Sheet "About" (sheet1): remove the closing paragraph (rows 15-18) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Rows("15:18").Select() | Out-Null
$wsAbout.Rows("15:18").Delete() | Out-Null

# --- Sheet "BAEPAbCiPC" (sheet2): flip several fuel flags from 0 to 1 ---
$wsData = $wb.Worksheets.Item("BAEPAbCiPC")
$wsData.Activate() | Out-Null

$wsData.Range("B3:B4").Value = 1
$wsData.Range("B9:B14").Value = 1
$wsData.Range("B17:B20").Value = 1

# Reflect the author's final selection / scroll position on this sheet
$wsData.Range("B17:B20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# The "About" sheet remains the tab that is selected/active when the
# workbook is saved (matches the original file's tabSelected state).
$wsAbout.Activate() | Out-Null
